$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Dll1"
$row2[0,2] = "Notch3"
$row2[0,3] = "ECs"
$row2[0,4] = [double]"3"
$row2[0,5] = [double]"1"
$row2[0,6] = [double]"16.39002933333333"
$row2[0,7] = [double]"49.170088"
$row2[0,8] = [double]"0.5551882184054378"
$row2[0,9] = [double]"0.5551882184054378"
$row2[0,10] = [double]"3"
$row2[0,11] = [double]"1"
$row2[0,12] = [double]"3.151158666666667"
$row2[0,13] = [double]"9.453476"
$row2[0,14] = [double]"0.03114707555614071"
$row2[0,15] = [double]"0.03114707555614071"
$row2[0,16] = [double]"51.64758298065423"
$row2[0,17] = [double]"464.828246825888"
$row2[0,18] = [double]"0.01729248938655332"
$row2[0,19] = [double]"0.01729248938655332"
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Dll1"
$row3[0,2] = "Notch3"
$row3[0,3] = "FAPs"
$row3[0,4] = [double]"3"
$row3[0,5] = [double]"1"
$row3[0,6] = [double]"16.39002933333333"
$row3[0,7] = [double]"49.170088"
$row3[0,8] = [double]"0.5551882184054378"
$row3[0,9] = [double]"0.5551882184054378"
$row3[0,10] = [double]"3"
$row3[0,11] = [double]"1"
$row3[0,12] = [double]"5.038243666666667"
$row3[0,13] = [double]"15.114731"
$row3[0,14] = [double]"0.04979963650066307"
$row3[0,15] = [double]"0.04979963650066306"
$row3[0,16] = [double]"82.57696148514756"
$row3[0,17] = [double]"743.1926533663282"
$row3[0,18] = [double]"0.02764817146604154"
$row3[0,19] = [double]"0.02764817146604153"
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "ECs"
$row4[0,1] = "Dll1"
$row4[0,2] = "Notch3"
$row4[0,3] = "MuSCs"
$row4[0,4] = [double]"3"
$row4[0,5] = [double]"1"
$row4[0,6] = [double]"16.39002933333333"
$row4[0,7] = [double]"49.170088"
$row4[0,8] = [double]"0.5551882184054378"
$row4[0,9] = [double]"0.5551882184054378"
$row4[0,10] = [double]"3"
$row4[0,11] = [double]"1"
$row4[0,12] = [double]"92.911639"
$row4[0,13] = [double]"278.734917"
$row4[0,14] = [double]"0.9183688116343246"
$row4[0,15] = [double]"0.9183688116343246"
$row4[0,16] = [double]"1522.824488618077"
$row4[0,17] = [double]"13705.4203975627"
$row4[0,18] = [double]"0.5098675443703797"
$row4[0,19] = [double]"0.5098675443703797"
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = "ECs"
$row5[0,1] = "Dll1"
$row5[0,2] = "Notch3"
$row5[0,3] = "Resolving-Mac"
$row5[0,4] = [double]"3"
$row5[0,5] = [double]"1"
$row5[0,6] = [double]"16.39002933333333"
$row5[0,7] = [double]"49.170088"
$row5[0,8] = [double]"0.5551882184054378"
$row5[0,9] = [double]"0.5551882184054378"
$row5[0,10] = [double]"2"
$row5[0,11] = [double]"0.6666666666666666"
$row5[0,12] = [double]"0.06924866666666667"
$row5[0,13] = [double]"0.207746"
$row5[0,14] = [double]"0.0006844763088715736"
$row5[0,15] = [double]"0.0006844763088715734"
$row5[0,16] = [double]"1.134987677960889"
$row5[0,17] = [double]"10.214889101648"
$row5[0,18] = [double]"0.0003800131824631391"
$row5[0,19] = [double]"0.000380013182463139"
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = "FAPs"
$row6[0,1] = "Dll1"
$row6[0,2] = "Notch3"
$row6[0,3] = "ECs"
$row6[0,4] = [double]"2"
$row6[0,5] = [double]"0.6666666666666666"
$row6[0,6] = [double]"0.2120556666666667"
$row6[0,7] = [double]"0.636167"
$row6[0,8] = [double]"0.007183074867352934"
$row6[0,9] = [double]"0.007183074867352935"
$row6[0,10] = [double]"3"
$row6[0,11] = [double]"1"
$row6[0,12] = [double]"3.151158666666667"
$row6[0,13] = [double]"9.453476"
$row6[0,14] = [double]"0.03114707555614071"
$row6[0,15] = [double]"0.03114707555614071"
$row6[0,16] = [double]"0.6682210518324445"
$row6[0,17] = [double]"6.013989466492"
$row6[0,18] = [double]"0.0002237317756188572"
$row6[0,19] = [double]"0.0002237317756188573"
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = "FAPs"
$row7[0,1] = "Dll1"
$row7[0,2] = "Notch3"
$row7[0,3] = "FAPs"
$row7[0,4] = [double]"2"
$row7[0,5] = [double]"0.6666666666666666"
$row7[0,6] = [double]"0.2120556666666667"
$row7[0,7] = [double]"0.636167"
$row7[0,8] = [double]"0.007183074867352934"
$row7[0,9] = [double]"0.007183074867352935"
$row7[0,10] = [double]"3"
$row7[0,11] = [double]"1"
$row7[0,12] = [double]"5.038243666666667"
$row7[0,13] = [double]"15.114731"
$row7[0,14] = [double]"0.04979963650066307"
$row7[0,15] = [double]"0.04979963650066306"
$row7[0,16] = [double]"1.068388119564111"
$row7[0,17] = [double]"9.615493076077001"
$row7[0,18] = [double]"0.0003577145173512247"
$row7[0,19] = [double]"0.0003577145173512247"
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,20
$row8[0,0] = "FAPs"
$row8[0,1] = "Dll1"
$row8[0,2] = "Notch3"
$row8[0,3] = "MuSCs"
$row8[0,4] = [double]"2"
$row8[0,5] = [double]"0.6666666666666666"
$row8[0,6] = [double]"0.2120556666666667"
$row8[0,7] = [double]"0.636167"
$row8[0,8] = [double]"0.007183074867352934"
$row8[0,9] = [double]"0.007183074867352935"
$row8[0,10] = [double]"3"
$row8[0,11] = [double]"1"
$row8[0,12] = [double]"92.911639"
$row8[0,13] = [double]"278.734917"
$row8[0,14] = [double]"0.9183688116343246"
$row8[0,15] = [double]"0.9183688116343246"
$row8[0,16] = [double]"19.70243954923767"
$row8[0,17] = [double]"177.321955943139"
$row8[0,18] = [double]"0.006596711929811297"
$row8[0,19] = [double]"0.006596711929811299"
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,20
$row9[0,0] = "FAPs"
$row9[0,1] = "Dll1"
$row9[0,2] = "Notch3"
$row9[0,3] = "Resolving-Mac"
$row9[0,4] = [double]"2"
$row9[0,5] = [double]"0.6666666666666666"
$row9[0,6] = [double]"0.2120556666666667"
$row9[0,7] = [double]"0.636167"
$row9[0,8] = [double]"0.007183074867352934"
$row9[0,9] = [double]"0.007183074867352935"
$row9[0,10] = [double]"2"
$row9[0,11] = [double]"0.6666666666666666"
$row9[0,12] = [double]"0.06924866666666667"
$row9[0,13] = [double]"0.207746"
$row9[0,14] = [double]"0.0006844763088715736"
$row9[0,15] = [double]"0.0006844763088715734"
$row9[0,16] = [double]"0.01468457217577778"
$row9[0,17] = [double]"0.132161149582"
$row9[0,18] = [double]"4.916644571553904e-06"
$row9[0,19] = [double]"4.916644571553904e-06"
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,20
$row10[0,0] = "MuSCs"
$row10[0,1] = "Dll1"
$row10[0,2] = "Notch3"
$row10[0,3] = "ECs"
$row10[0,4] = [double]"3"
$row10[0,5] = [double]"1"
$row10[0,6] = [double]"12.89091033333333"
$row10[0,7] = [double]"38.672731"
$row10[0,8] = [double]"0.4366606914505164"
$row10[0,9] = [double]"0.4366606914505165"
$row10[0,10] = [double]"3"
$row10[0,11] = [double]"1"
$row10[0,12] = [double]"3.151158666666667"
$row10[0,13] = [double]"9.453476"
$row10[0,14] = [double]"0.03114707555614071"
$row10[0,15] = [double]"0.03114707555614071"
$row10[0,16] = [double]"40.62130381810622"
$row10[0,17] = [double]"365.591734362956"
$row10[0,18] = [double]"0.01360070354900588"
$row10[0,19] = [double]"0.01360070354900588"
$ws.Range("A10:T10").Value = $row10

$row11 = New-Object 'object[,]' 1,20
$row11[0,0] = "MuSCs"
$row11[0,1] = "Dll1"
$row11[0,2] = "Notch3"
$row11[0,3] = "FAPs"
$row11[0,4] = [double]"3"
$row11[0,5] = [double]"1"
$row11[0,6] = [double]"12.89091033333333"
$row11[0,7] = [double]"38.672731"
$row11[0,8] = [double]"0.4366606914505164"
$row11[0,9] = [double]"0.4366606914505165"
$row11[0,10] = [double]"3"
$row11[0,11] = [double]"1"
$row11[0,12] = [double]"5.038243666666667"
$row11[0,13] = [double]"15.114731"
$row11[0,14] = [double]"0.04979963650066307"
$row11[0,15] = [double]"0.04979963650066306"
$row11[0,16] = [double]"64.94754734448456"
$row11[0,17] = [double]"584.5279261003611"
$row11[0,18] = [double]"0.02174554370836391"
$row11[0,19] = [double]"0.02174554370836391"
$ws.Range("A11:T11").Value = $row11

$row12 = New-Object 'object[,]' 1,20
$row12[0,0] = "MuSCs"
$row12[0,1] = "Dll1"
$row12[0,2] = "Notch3"
$row12[0,3] = "MuSCs"
$row12[0,4] = [double]"3"
$row12[0,5] = [double]"1"
$row12[0,6] = [double]"12.89091033333333"
$row12[0,7] = [double]"38.672731"
$row12[0,8] = [double]"0.4366606914505164"
$row12[0,9] = [double]"0.4366606914505165"
$row12[0,10] = [double]"3"
$row12[0,11] = [double]"1"
$row12[0,12] = [double]"92.911639"
$row12[0,13] = [double]"278.734917"
$row12[0,14] = [double]"0.9183688116343246"
$row12[0,15] = [double]"0.9183688116343246"
$row12[0,16] = [double]"1197.715607272036"
$row12[0,17] = [double]"10779.44046544833"
$row12[0,18] = [double]"0.4010155602948332"
$row12[0,19] = [double]"0.4010155602948333"
$ws.Range("A12:T12").Value = $row12

$row13 = New-Object 'object[,]' 1,20
$row13[0,0] = "MuSCs"
$row13[0,1] = "Dll1"
$row13[0,2] = "Notch3"
$row13[0,3] = "Resolving-Mac"
$row13[0,4] = [double]"3"
$row13[0,5] = [double]"1"
$row13[0,6] = [double]"12.89091033333333"
$row13[0,7] = [double]"38.672731"
$row13[0,8] = [double]"0.4366606914505164"
$row13[0,9] = [double]"0.4366606914505165"
$row13[0,10] = [double]"2"
$row13[0,11] = [double]"0.6666666666666666"
$row13[0,12] = [double]"0.06924866666666667"
$row13[0,13] = [double]"0.207746"
$row13[0,14] = [double]"0.0006844763088715736"
$row13[0,15] = [double]"0.0006844763088715734"
$row13[0,16] = [double]"0.8926783527028889"
$row13[0,17] = [double]"8.034105174326"
$row13[0,18] = [double]"0.0002988838983133585"
$row13[0,19] = [double]"0.0002988838983133585"
$ws.Range("A13:T13").Value = $row13

$row14 = New-Object 'object[,]' 1,20
$row14[0,0] = "Resolving-Mac"
$row14[0,1] = "Dll1"
$row14[0,2] = "Notch3"
$row14[0,3] = "ECs"
$row14[0,4] = [double]"1"
$row14[0,5] = [double]"0.3333333333333333"
$row14[0,6] = [double]"0.02857733333333333"
$row14[0,7] = [double]"0.085732"
$row14[0,8] = [double]"0.0009680152766929151"
$row14[0,9] = [double]"0.0009680152766929153"
$row14[0,10] = [double]"3"
$row14[0,11] = [double]"1"
$row14[0,12] = [double]"3.151158666666667"
$row14[0,13] = [double]"9.453476"
$row14[0,14] = [double]"0.03114707555614071"
$row14[0,15] = [double]"0.03114707555614071"
$row14[0,16] = [double]"0.09005171160355556"
$row14[0,17] = [double]"0.810465404432"
$row14[0,18] = [double]"3.015084496265269e-05"
$row14[0,19] = [double]"3.015084496265269e-05"
$ws.Range("A14:T14").Value = $row14

$row15 = New-Object 'object[,]' 1,20
$row15[0,0] = "Resolving-Mac"
$row15[0,1] = "Dll1"
$row15[0,2] = "Notch3"
$row15[0,3] = "FAPs"
$row15[0,4] = [double]"1"
$row15[0,5] = [double]"0.3333333333333333"
$row15[0,6] = [double]"0.02857733333333333"
$row15[0,7] = [double]"0.085732"
$row15[0,8] = [double]"0.0009680152766929151"
$row15[0,9] = [double]"0.0009680152766929153"
$row15[0,10] = [double]"3"
$row15[0,11] = [double]"1"
$row15[0,12] = [double]"5.038243666666667"
$row15[0,13] = [double]"15.114731"
$row15[0,14] = [double]"0.04979963650066307"
$row15[0,15] = [double]"0.04979963650066306"
$row15[0,16] = [double]"0.1439795686768889"
$row15[0,17] = [double]"1.295816118092"
$row15[0,18] = [double]"4.820680890639595e-05"
$row15[0,19] = [double]"4.820680890639596e-05"
$ws.Range("A15:T15").Value = $row15

$row16 = New-Object 'object[,]' 1,20
$row16[0,0] = "Resolving-Mac"
$row16[0,1] = "Dll1"
$row16[0,2] = "Notch3"
$row16[0,3] = "MuSCs"
$row16[0,4] = [double]"1"
$row16[0,5] = [double]"0.3333333333333333"
$row16[0,6] = [double]"0.02857733333333333"
$row16[0,7] = [double]"0.085732"
$row16[0,8] = [double]"0.0009680152766929151"
$row16[0,9] = [double]"0.0009680152766929153"
$row16[0,10] = [double]"3"
$row16[0,11] = [double]"1"
$row16[0,12] = [double]"92.911639"
$row16[0,13] = [double]"278.734917"
$row16[0,14] = [double]"0.9183688116343246"
$row16[0,15] = [double]"0.9183688116343246"
$row16[0,16] = [double]"2.655166878249333"
$row16[0,17] = [double]"23.896501904244"
$row16[0,18] = [double]"0.0008889950393003444"
$row16[0,19] = [double]"0.0008889950393003446"
$ws.Range("A16:T16").Value = $row16

$row17 = New-Object 'object[,]' 1,20
$row17[0,0] = "Resolving-Mac"
$row17[0,1] = "Dll1"
$row17[0,2] = "Notch3"
$row17[0,3] = "Resolving-Mac"
$row17[0,4] = [double]"1"
$row17[0,5] = [double]"0.3333333333333333"
$row17[0,6] = [double]"0.02857733333333333"
$row17[0,7] = [double]"0.085732"
$row17[0,8] = [double]"0.0009680152766929151"
$row17[0,9] = [double]"0.0009680152766929153"
$row17[0,10] = [double]"2"
$row17[0,11] = [double]"0.6666666666666666"
$row17[0,12] = [double]"0.06924866666666667"
$row17[0,13] = [double]"0.207746"
$row17[0,14] = [double]"0.0006844763088715736"
$row17[0,15] = [double]"0.0006844763088715734"
$row17[0,16] = [double]"0.001978942230222222"
$row17[0,17] = [double]"0.017810480072"
$row17[0,18] = [double]"6.625835235220615e-07"
$row17[0,19] = [double]"6.625835235220615e-07"
$ws.Range("A17:T17").Value = $row17
